# Update pivots_long sheet: refresh occurrence counts/percentages for
# EURUSD/GBPUSD/USDJPY tables, and shrink the XAUUSD table's year range
# from 2016-2022 (7 years) down to 2017-2022 (6 years, matching the other
# instrument tables), shifting its total/percentage columns left by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pivots_long")

# ---------------------------------------------------------------
# EURUSD table (rows 4-7): year columns B:H (2016-2022) unchanged,
# total = I, percentage = J
# ---------------------------------------------------------------
$eurusd = @{
    4 = @(70, 81, 106, 61, 114, 54, 123, 609, 33.36986301369863)
    5 = @(25, 32, 38, 21, 57, 38, 82, 293, 16.05479452054795)
    6 = @(133, 117, 96, 153, 73, 102, 36, 710, 38.9041095890411)
    7 = @(32, 29, 19, 24, 18, 69, 22, 213, 11.67123287671233)
}
foreach ($row in $eurusd.Keys) {
    $vals = $eurusd[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($row, 2 + $i).Value = $vals[$i]
    }
}

# ---------------------------------------------------------------
# GBPUSD table (rows 11-14): year columns B:G (2017-2022) unchanged,
# total = H, percentage = I
# ---------------------------------------------------------------
$gbpusd = @{
    11 = @(80, 99, 51, 109, 54, 129, 522, 33.35463258785942)
    12 = @(19, 42, 63, 58, 20, 70, 272, 17.38019169329073)
    13 = @(136, 93, 95, 77, 162, 53, 616, 39.36102236421725)
    14 = @(24, 25, 50, 18, 27, 11, 155, 9.904153354632587)
}
foreach ($row in $gbpusd.Keys) {
    $vals = $gbpusd[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($row, 2 + $i).Value = $vals[$i]
    }
}

# ---------------------------------------------------------------
# USDJPY table (rows 18-21): year columns B:G (2017-2022) unchanged,
# total = H, percentage = I
# ---------------------------------------------------------------
$usdjpy = @{
    18 = @(57, 66, 75, 81, 93, 124, 496, 31.69329073482428)
    19 = @(27, 33, 35, 57, 40, 79, 271, 17.31629392971246)
    20 = @(129, 116, 123, 99, 102, 40, 609, 38.91373801916933)
    21 = @(46, 44, 26, 25, 28, 20, 189, 12.07667731629393)
}
foreach ($row in $usdjpy.Keys) {
    $vals = $usdjpy[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($row, 2 + $i).Value = $vals[$i]
    }
}

# ---------------------------------------------------------------
# XAUUSD table: header row 23 shrinks from 7 years (B:H = 2016-2022,
# total = I, percentage = J) to 6 years (B:G = 2017-2022, total = H,
# percentage = I), matching the GBPUSD/USDJPY layout. Clear the
# now-unused column J first, then rewrite the header + data rows
# (25-28) with the shifted values.
# ---------------------------------------------------------------
$ws.Range("B23:J23").ClearContents()
$ws.Range("B25:J28").ClearContents()

$years = @(2017, 2018, 2019, 2020, 2021, 2022)
for ($i = 0; $i -lt $years.Length; $i++) {
    $ws.Cells.Item(23, 2 + $i).Value = $years[$i]
}
$ws.Cells.Item(23, 8).Value = "total_count_of_occurrences"
$ws.Cells.Item(23, 9).Value = "percentage_of_occurrences"

$xauusd = @{
    25 = @(72, 101, 90, 125, 66, 95, 549, 35.46511627906977)
    26 = @(21, 19, 51, 38, 40, 58, 227, 14.6640826873385)
    27 = @(114, 119, 83, 77, 104, 68, 565, 36.49870801033592)
    28 = @(50, 19, 34, 19, 48, 37, 207, 13.37209302325581)
}
foreach ($row in $xauusd.Keys) {
    $vals = $xauusd[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($row, 2 + $i).Value = $vals[$i]
    }
}
